$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 51, pushing existing rows 51-55 down to 53-57
$ws.Range("A51:A52").EntireRow.Insert()

# Fill in the two new rows (51 and 52) with the new weekly price data
# Row 51: Packham's Triumph
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 45021
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100104
$ws.Range("H51").Value = "Frutos de pepita"
$ws.Range("I51").Value = 100104005
$ws.Range("J51").Value = "Pera"
$ws.Range("K51").Value = "Packham's Triumph"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 300
$ws.Range("N51").Value = 19000
$ws.Range("O51").Value = 20000
$ws.Range("P51").Value = 19500
$ws.Range("Q51").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 1083
$ws.Range("T51").Value = 18

# Row 52: Winter Nelis
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 45021
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100104
$ws.Range("H52").Value = "Frutos de pepita"
$ws.Range("I52").Value = 100104005
$ws.Range("J52").Value = "Pera"
$ws.Range("K52").Value = "Winter Nelis"
$ws.Range("L52").Value = "Segunda"
$ws.Range("M52").Value = 300
$ws.Range("N52").Value = 19000
$ws.Range("O52").Value = 20000
$ws.Range("P52").Value = 19500
$ws.Range("Q52").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R52").Value = "Región de O'Higgins"
$ws.Range("S52").Value = 1083
$ws.Range("T52").Value = 18

# Make sure the date cells keep/get the date number format used elsewhere in column D
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
